$d = $word.ActiveDocument

function Get-ParagraphIndexByStart($doc, $startPos) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs($i).Range.Start -eq $startPos) {
            return $i
        }
    }
    return -1
}

# ------------------------------------------------------------------
# 1) Insert a new bullet right after the paragraph that talks about
#    elevator details being stored in the database (i.e. right before
#    the "Database Design" Heading2), reusing the same list numbering
#    (numId 1 / ListParagraph style) as its neighbour.
# ------------------------------------------------------------------
$anchor1 = $d.Content.Duplicate
$anchor1.Find.Execute("Elevator details will be stored in the database", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$anchorPara1 = $anchor1.Paragraphs(1)
$anchorIdx1 = Get-ParagraphIndexByStart $d $anchorPara1.Range.Start

$anchorPara1.Range.InsertParagraphAfter() | Out-Null
$newIdx1 = $anchorIdx1 + 1

$newRange1 = $d.Paragraphs($newIdx1).Range
$newRange1.InsertAfter("Although console " + [char]0x201C + "ReadLine" + [char]0x201D + " operations aren" + [char]0x2019 + "t strictly ")

$tailRange1 = $d.Paragraphs($newIdx1).Range
$insertionPoint1 = $d.Range($tailRange1.End - 1, $tailRange1.End - 1)
$insertionPoint1.InsertAfter("async, operations in the master elevator control as async to allow for future extensions with different UI types, or multiple elevator calls")

# ------------------------------------------------------------------
# 2) Insert a new bullet right after the paragraph about dictionaries
#    of phrases / localization, at the very end of the document,
#    reusing the same list numbering (numId 4 / ListParagraph style).
# ------------------------------------------------------------------
$anchor2 = $d.Content.Duplicate
$anchor2.Find.Execute("Dictionaries of phrases can be added for different languages", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$anchorPara2 = $anchor2.Paragraphs(1)
$anchorIdx2 = Get-ParagraphIndexByStart $d $anchorPara2.Range.Start

$anchorPara2.Range.InsertParagraphAfter() | Out-Null
$newIdx2 = $anchorIdx2 + 1

$newRange2 = $d.Paragraphs($newIdx2).Range
$newRange2.InsertAfter("A logging mechanism that can be switched between verbose, dev and prod")
